$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New day of data (2025-09-22, serial 45922) for both sites, appended
# after the existing last row (43).
$newRows = @(
    @{ Row = 44; Date = 45922; Site = "四方坪站"; C = 8088.95; D = 6532.5;  E = 2748;    F = 352 },
    @{ Row = 45; Date = 45922; Site = "高岭站";   C = 3967.86; D = 3125;    E = 1040.42; F = 151 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    # Leave number formatting alone: the columns already carry the correct
    # style (via the sheet's <cols> definitions), and new cells in those
    # columns pick it up automatically, just like Excel does when you type
    # values into the row right below existing data.
    $ws.Cells.Item($row, 1).Value = $r.Date
    $ws.Cells.Item($row, 2).Value = $r.Site
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
}

# Update the view to match the author's final scroll/selection position.
$excel.ActiveWindow.ScrollRow = 34
$ws.Range("H43").Select()
